$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" -> "_FV2404", "_new" -> "_FV2410" -----------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $current = [string]$cell.Value2
    if ($current -ne $null) {
        if ($current.EndsWith("_old")) {
            $cell.Value = $current.Substring(0, $current.Length - 4) + "_FV2404"
        } elseif ($current.EndsWith("_new")) {
            $cell.Value = $current.Substring(0, $current.Length - 4) + "_FV2410"
        }
    }
}

# --- 2. Turn the used range into an Excel Table named "Table1" ---------------
$tableRange = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split after row 1, top-left cell A2) ----------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header renaming, table creation and freeze panes applied."
